$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 48.091872
$ws.Cells.Item(2, 8).Value = 144.275616
$ws.Cells.Item(2, 9).Value = 0.421093842675958
$ws.Cells.Item(2, 10).Value = 0.423782205092405
$ws.Cells.Item(2, 13).Value = 2.621797333333333
$ws.Cells.Item(2, 14).Value = 7.865392
$ws.Cells.Item(2, 15).Value = 0.07867217155043885
$ws.Cells.Item(2, 16).Value = 0.07906089226781998
$ws.Cells.Item(2, 17).Value = 126.087141764608
$ws.Cells.Item(2, 18).Value = 1134.784275881472
$ws.Cells.Item(2, 19).Value = 0.03312836702983647
$ws.Cells.Item(2, 20).Value = 0.03350459926182982
# Row 3
$ws.Cells.Item(3, 7).Value = 48.091872
$ws.Cells.Item(3, 8).Value = 144.275616
$ws.Cells.Item(3, 9).Value = 0.421093842675958
$ws.Cells.Item(3, 10).Value = 0.423782205092405
$ws.Cells.Item(3, 15).Value = 0.148308476825081
$ws.Cells.Item(3, 16).Value = 0.1490412718702539
$ws.Cells.Item(3, 17).Value = 237.692586511968
$ws.Cells.Item(3, 18).Value = 2139.233278607712
$ws.Cells.Item(3, 19).Value = 0.06245178640769163
$ws.Cells.Item(3, 20).Value = 0.06316103884295284
# Row 4
$ws.Cells.Item(4, 7).Value = 48.091872
$ws.Cells.Item(4, 8).Value = 144.275616
$ws.Cells.Item(4, 9).Value = 0.421093842675958
$ws.Cells.Item(4, 10).Value = 0.423782205092405
$ws.Cells.Item(4, 13).Value = 11.09754033333333
$ws.Cells.Item(4, 14).Value = 33.292621
$ws.Cells.Item(4, 15).Value = 0.3330034651388949
$ws.Cells.Item(4, 16).Value = 0.3346488416844782
$ws.Cells.Item(4, 17).Value = 533.701489225504
$ws.Cells.Item(4, 18).Value = 4803.313403029536
$ws.Cells.Item(4, 19).Value = 0.1402257087597467
$ws.Cells.Item(4, 20).Value = 0.1418182240606673
# Row 5
$ws.Cells.Item(5, 7).Value = 48.091872
$ws.Cells.Item(5, 8).Value = 144.275616
$ws.Cells.Item(5, 9).Value = 0.421093842675958
$ws.Cells.Item(5, 10).Value = 0.423782205092405
$ws.Cells.Item(5, 13).Value = 0.4915585
$ws.Cells.Item(5, 14).Value = 0.983117
$ws.Cells.Item(5, 15).Value = 0.01475017696730553
$ws.Cells.Item(5, 16).Value = 0.009882038584175128
$ws.Cells.Item(5, 17).Value = 23.639968462512
$ws.Cells.Item(5, 18).Value = 141.839810775072
$ws.Cells.Item(5, 19).Value = 0.006211208699313094
$ws.Cells.Item(5, 20).Value = 0.004187832102009964
# Row 6
$ws.Cells.Item(6, 7).Value = 48.091872
$ws.Cells.Item(6, 8).Value = 144.275616
$ws.Cells.Item(6, 9).Value = 0.421093842675958
$ws.Cells.Item(6, 10).Value = 0.423782205092405
$ws.Cells.Item(6, 13).Value = 14.172235
$ws.Cells.Item(6, 14).Value = 42.516705
$ws.Cells.Item(6, 15).Value = 0.4252657095182797
$ws.Cells.Item(6, 16).Value = 0.4273669555932728
$ws.Cells.Item(6, 17).Value = 681.56931157392
$ws.Cells.Item(6, 18).Value = 6134.12380416528
$ws.Cells.Item(6, 19).Value = 0.1790767717793701
$ws.Cells.Item(6, 20).Value = 0.1811105108249451
# Row 7
$ws.Cells.Item(7, 9).Value = 0.1230362686979479
$ws.Cells.Item(7, 10).Value = 0.1238217612582891
$ws.Cells.Item(7, 13).Value = 2.621797333333333
$ws.Cells.Item(7, 14).Value = 7.865392
$ws.Cells.Item(7, 15).Value = 0.07867217155043885
$ws.Cells.Item(7, 16).Value = 0.07906089226781998
$ws.Cells.Item(7, 17).Value = 36.84046139198578
$ws.Cells.Item(7, 18).Value = 331.564152527872
$ws.Cells.Item(7, 19).Value = 0.00967953043793085
$ws.Cells.Item(7, 20).Value = 0.009789458927253318
# Row 8
$ws.Cells.Item(8, 9).Value = 0.1230362686979479
$ws.Cells.Item(8, 10).Value = 0.1238217612582891
$ws.Cells.Item(8, 15).Value = 0.148308476825081
$ws.Cells.Item(8, 16).Value = 0.1490412718702539
$ws.Cells.Item(8, 19).Value = 0.01824732160483405
$ws.Cells.Item(8, 20).Value = 0.01845455278315033
# Row 9
$ws.Cells.Item(9, 9).Value = 0.1230362686979479
$ws.Cells.Item(9, 10).Value = 0.1238217612582891
$ws.Cells.Item(9, 13).Value = 11.09754033333333
$ws.Cells.Item(9, 14).Value = 33.292621
$ws.Cells.Item(9, 15).Value = 0.3330034651388949
$ws.Cells.Item(9, 16).Value = 0.3346488416844782
$ws.Cells.Item(9, 17).Value = 155.9382569347484
$ws.Cells.Item(9, 18).Value = 1403.444312412736
$ws.Cells.Item(9, 19).Value = 0.04097150381417681
$ws.Cells.Item(9, 20).Value = 0.04143680898041843
# Row 10
$ws.Cells.Item(10, 9).Value = 0.1230362686979479
$ws.Cells.Item(10, 10).Value = 0.1238217612582891
$ws.Cells.Item(10, 13).Value = 0.4915585
$ws.Cells.Item(10, 14).Value = 0.983117
$ws.Cells.Item(10, 15).Value = 0.01475017696730553
$ws.Cells.Item(10, 16).Value = 0.009882038584175128
$ws.Cells.Item(10, 17).Value = 6.907186040245334
$ws.Cells.Item(10, 18).Value = 41.443116241472
$ws.Cells.Item(10, 19).Value = 0.001814806736691686
$ws.Cells.Item(10, 20).Value = 0.001223611422314934
# Row 11
$ws.Cells.Item(11, 9).Value = 0.1230362686979479
$ws.Cells.Item(11, 10).Value = 0.1238217612582891
$ws.Cells.Item(11, 13).Value = 14.172235
$ws.Cells.Item(11, 14).Value = 42.516705
$ws.Cells.Item(11, 15).Value = 0.4252657095182797
$ws.Cells.Item(11, 16).Value = 0.4273669555932728
$ws.Cells.Item(11, 17).Value = 199.1426529112533
$ws.Cells.Item(11, 18).Value = 1792.28387620128
$ws.Cells.Item(11, 19).Value = 0.05232310610431454
$ws.Cells.Item(11, 20).Value = 0.05291732914515206
# Row 12
$ws.Cells.Item(12, 7).Value = 21.412221
$ws.Cells.Item(12, 8).Value = 64.23666299999999
$ws.Cells.Item(12, 9).Value = 0.1874860355013181
$ws.Cells.Item(12, 10).Value = 0.1886829905749125
$ws.Cells.Item(12, 13).Value = 2.621797333333333
$ws.Cells.Item(12, 14).Value = 7.865392
$ws.Cells.Item(12, 15).Value = 0.07867217155043885
$ws.Cells.Item(12, 16).Value = 0.07906089226781998
$ws.Cells.Item(12, 17).Value = 56.138503918544
$ws.Cells.Item(12, 18).Value = 505.2465352668959
$ws.Cells.Item(12, 19).Value = 0.01474993354827136
$ws.Cells.Item(12, 20).Value = 0.01491744559061325
# Row 13
$ws.Cells.Item(13, 7).Value = 21.412221
$ws.Cells.Item(13, 8).Value = 64.23666299999999
$ws.Cells.Item(13, 9).Value = 0.1874860355013181
$ws.Cells.Item(13, 10).Value = 0.1886829905749125
$ws.Cells.Item(13, 15).Value = 0.148308476825081
$ws.Cells.Item(13, 16).Value = 0.1490412718702539
$ws.Cells.Item(13, 17).Value = 105.829238513649
$ws.Cells.Item(13, 18).Value = 952.463146622841
$ws.Cells.Item(13, 19).Value = 0.02780576835117355
$ws.Cells.Item(13, 20).Value = 0.02812155289556809
# Row 14
$ws.Cells.Item(14, 7).Value = 21.412221
$ws.Cells.Item(14, 8).Value = 64.23666299999999
$ws.Cells.Item(14, 9).Value = 0.1874860355013181
$ws.Cells.Item(14, 10).Value = 0.1886829905749125
$ws.Cells.Item(14, 13).Value = 11.09754033333333
$ws.Cells.Item(14, 14).Value = 33.292621
$ws.Cells.Item(14, 15).Value = 0.3330034651388949
$ws.Cells.Item(14, 16).Value = 0.3346488416844782
$ws.Cells.Item(14, 17).Value = 237.622986173747
$ws.Cells.Item(14, 18).Value = 2138.606875563722
$ws.Cells.Item(14, 19).Value = 0.06243349948709278
$ws.Cells.Item(14, 20).Value = 0.06314254424145778
# Row 15
$ws.Cells.Item(15, 7).Value = 21.412221
$ws.Cells.Item(15, 8).Value = 64.23666299999999
$ws.Cells.Item(15, 9).Value = 0.1874860355013181
$ws.Cells.Item(15, 10).Value = 0.1886829905749125
$ws.Cells.Item(15, 13).Value = 0.4915585
$ws.Cells.Item(15, 14).Value = 0.983117
$ws.Cells.Item(15, 15).Value = 0.01475017696730553
$ws.Cells.Item(15, 16).Value = 0.009882038584175128
$ws.Cells.Item(15, 17).Value = 10.5253592364285
$ws.Cells.Item(15, 18).Value = 63.15215541857099
$ws.Cells.Item(15, 19).Value = 0.002765452202542968
$ws.Cells.Item(15, 20).Value = 0.001864572593038837
# Row 16
$ws.Cells.Item(16, 7).Value = 21.412221
$ws.Cells.Item(16, 8).Value = 64.23666299999999
$ws.Cells.Item(16, 9).Value = 0.1874860355013181
$ws.Cells.Item(16, 10).Value = 0.1886829905749125
$ws.Cells.Item(16, 13).Value = 14.172235
$ws.Cells.Item(16, 14).Value = 42.516705
$ws.Cells.Item(16, 15).Value = 0.4252657095182797
$ws.Cells.Item(16, 16).Value = 0.4273669555932728
$ws.Cells.Item(16, 17).Value = 303.459027883935
$ws.Cells.Item(16, 18).Value = 2731.131250955415
$ws.Cells.Item(16, 19).Value = 0.0797313819122374
$ws.Cells.Item(16, 20).Value = 0.08063687525423455
# Row 17
$ws.Cells.Item(17, 7).Value = 2.1734975
$ws.Cells.Item(17, 8).Value = 4.346995
$ws.Cells.Item(17, 9).Value = 0.01903120789977957
$ws.Cells.Item(17, 10).Value = 0.012768471746644
$ws.Cells.Item(17, 13).Value = 2.621797333333333
$ws.Cells.Item(17, 14).Value = 7.865392
$ws.Cells.Item(17, 15).Value = 0.07867217155043885
$ws.Cells.Item(17, 16).Value = 0.07906089226781998
$ws.Cells.Item(17, 17).Value = 5.698469949506666
$ws.Cells.Item(17, 18).Value = 34.19081969704
$ws.Cells.Item(17, 19).Value = 0.001497226452703525
$ws.Cells.Item(17, 20).Value = 0.001009486769186124
# Row 18
$ws.Cells.Item(18, 7).Value = 2.1734975
$ws.Cells.Item(18, 8).Value = 4.346995
$ws.Cells.Item(18, 9).Value = 0.01903120789977957
$ws.Cells.Item(18, 10).Value = 0.012768471746644
$ws.Cells.Item(18, 15).Value = 0.148308476825081
$ws.Cells.Item(18, 16).Value = 0.1490412718702539
$ws.Cells.Item(18, 17).Value = 10.7424440153275
$ws.Cells.Item(18, 18).Value = 64.454664091965
$ws.Cells.Item(18, 19).Value = 0.002822489455757758
$ws.Cells.Item(18, 20).Value = 0.001903029268959224
# Row 19
$ws.Cells.Item(19, 7).Value = 2.1734975
$ws.Cells.Item(19, 8).Value = 4.346995
$ws.Cells.Item(19, 9).Value = 0.01903120789977957
$ws.Cells.Item(19, 10).Value = 0.012768471746644
$ws.Cells.Item(19, 13).Value = 11.09754033333333
$ws.Cells.Item(19, 14).Value = 33.292621
$ws.Cells.Item(19, 15).Value = 0.3330034651388949
$ws.Cells.Item(19, 16).Value = 0.3346488416844782
$ws.Cells.Item(19, 17).Value = 24.12047617064917
$ws.Cells.Item(19, 18).Value = 144.722857023895
$ws.Cells.Item(19, 19).Value = 0.006337458176405308
$ws.Cells.Item(19, 20).Value = 0.0042729542800954
# Row 20
$ws.Cells.Item(20, 7).Value = 2.1734975
$ws.Cells.Item(20, 8).Value = 4.346995
$ws.Cells.Item(20, 9).Value = 0.01903120789977957
$ws.Cells.Item(20, 10).Value = 0.012768471746644
$ws.Cells.Item(20, 13).Value = 0.4915585
$ws.Cells.Item(20, 14).Value = 0.983117
$ws.Cells.Item(20, 15).Value = 0.01475017696730553
$ws.Cells.Item(20, 16).Value = 0.009882038584175128
$ws.Cells.Item(20, 17).Value = 1.06840117085375
$ws.Cells.Item(20, 18).Value = 4.273604683415
$ws.Cells.Item(20, 19).Value = 0.0002807136844233317
$ws.Cells.Item(20, 20).Value = 0.000126178530461286
# Row 21
$ws.Cells.Item(21, 7).Value = 2.1734975
$ws.Cells.Item(21, 8).Value = 4.346995
$ws.Cells.Item(21, 9).Value = 0.01903120789977957
$ws.Cells.Item(21, 10).Value = 0.012768471746644
$ws.Cells.Item(21, 13).Value = 14.172235
$ws.Cells.Item(21, 14).Value = 42.516705
$ws.Cells.Item(21, 15).Value = 0.4252657095182797
$ws.Cells.Item(21, 16).Value = 0.4273669555932728
$ws.Cells.Item(21, 17).Value = 30.8033173419125
$ws.Cells.Item(21, 18).Value = 184.819904051475
$ws.Cells.Item(21, 19).Value = 0.00809332013048965
$ws.Cells.Item(21, 20).Value = 0.005456822897941964
# Row 22
$ws.Cells.Item(22, 7).Value = 28.477822
$ws.Cells.Item(22, 8).Value = 85.433466
$ws.Cells.Item(22, 9).Value = 0.2493526452249964
$ws.Cells.Item(22, 10).Value = 0.2509445713277496
$ws.Cells.Item(22, 13).Value = 2.621797333333333
$ws.Cells.Item(22, 14).Value = 7.865392
$ws.Cells.Item(22, 15).Value = 0.07867217155043885
$ws.Cells.Item(22, 16).Value = 0.07906089226781998
$ws.Cells.Item(22, 17).Value = 74.66307777874133
$ws.Cells.Item(22, 18).Value = 671.967700008672
$ws.Cells.Item(22, 19).Value = 0.01961711408169663
$ws.Cells.Item(22, 20).Value = 0.01983990171893748
# Row 23
$ws.Cells.Item(23, 7).Value = 28.477822
$ws.Cells.Item(23, 8).Value = 85.433466
$ws.Cells.Item(23, 9).Value = 0.2493526452249964
$ws.Cells.Item(23, 10).Value = 0.2509445713277496
$ws.Cells.Item(23, 15).Value = 0.148308476825081
$ws.Cells.Item(23, 16).Value = 0.1490412718702539
$ws.Cells.Item(23, 17).Value = 140.750752422518
$ws.Cells.Item(23, 18).Value = 1266.756771802662
$ws.Cells.Item(23, 19).Value = 0.03698111100562403
$ws.Cells.Item(23, 20).Value = 0.03740109807962344
# Row 24
$ws.Cells.Item(24, 7).Value = 28.477822
$ws.Cells.Item(24, 8).Value = 85.433466
$ws.Cells.Item(24, 9).Value = 0.2493526452249964
$ws.Cells.Item(24, 10).Value = 0.2509445713277496
$ws.Cells.Item(24, 13).Value = 11.09754033333333
$ws.Cells.Item(24, 14).Value = 33.292621
$ws.Cells.Item(24, 15).Value = 0.3330034651388949
$ws.Cells.Item(24, 16).Value = 0.3346488416844782
$ws.Cells.Item(24, 17).Value = 316.0337782504873
$ws.Cells.Item(24, 18).Value = 2844.304004254386
$ws.Cells.Item(24, 19).Value = 0.08303529490147331
$ws.Cells.Item(24, 20).Value = 0.0839783101218393
# Row 25
$ws.Cells.Item(25, 7).Value = 28.477822
$ws.Cells.Item(25, 8).Value = 85.433466
$ws.Cells.Item(25, 9).Value = 0.2493526452249964
$ws.Cells.Item(25, 10).Value = 0.2509445713277496
$ws.Cells.Item(25, 13).Value = 0.4915585
$ws.Cells.Item(25, 14).Value = 0.983117
$ws.Cells.Item(25, 15).Value = 0.01475017696730553
$ws.Cells.Item(25, 16).Value = 0.009882038584175128
$ws.Cells.Item(25, 17).Value = 13.998515465587
$ws.Cells.Item(25, 18).Value = 83.991092793522
$ws.Cells.Item(25, 19).Value = 0.003677995644334449
$ws.Cells.Item(25, 20).Value = 0.002479843936350108
# Row 26
$ws.Cells.Item(26, 7).Value = 28.477822
$ws.Cells.Item(26, 8).Value = 85.433466
$ws.Cells.Item(26, 9).Value = 0.2493526452249964
$ws.Cells.Item(26, 10).Value = 0.2509445713277496
$ws.Cells.Item(26, 13).Value = 14.172235
$ws.Cells.Item(26, 14).Value = 42.516705
$ws.Cells.Item(26, 15).Value = 0.4252657095182797
$ws.Cells.Item(26, 16).Value = 0.4273669555932728
$ws.Cells.Item(26, 17).Value = 403.59438567217
$ws.Cells.Item(26, 18).Value = 3632.34947104953
$ws.Cells.Item(26, 19).Value = 0.106041129591868
$ws.Cells.Item(26, 20).Value = 0.1072454174709992
